$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 23827156
$ws.Cells.Item(40, 10).Value = 71447580
$ws.Cells.Item(40, 12).Value = 71447580
$ws.Cells.Item(40, 14).Value = -71447930
$ws.Cells.Item(64, 8).Value = 25253332
$ws.Cells.Item(64, 10).Value = 50004296
$ws.Cells.Item(64, 12).Value = 50004296
$ws.Cells.Item(64, 14).Value = -50004792
$ws.Cells.Item(67, 8).Value = 25253332
$ws.Cells.Item(67, 10).Value = 50004296
$ws.Cells.Item(67, 12).Value = 50004296
$ws.Cells.Item(67, 14).Value = -50006012
$ws.Cells.Item(74, 9).Value = 5000
$ws.Cells.Item(74, 11).Value = 5000
$ws.Cells.Item(74, 13).Value = -4064
$ws.Cells.Item(77, 9).Value = 5000
$ws.Cells.Item(77, 11).Value = 25000
$ws.Cells.Item(77, 13).Value = -20320
$ws.Cells.Item(80, 8).Value = 793.75
$ws.Cells.Item(80, 10).Value = 675
$ws.Cells.Item(80, 12).Value = 2025
$ws.Cells.Item(80, 14).Value = -4021
$ws.Cells.Item(83, 8).Value = 793.75
$ws.Cells.Item(83, 10).Value = 675
$ws.Cells.Item(83, 12).Value = 6075
$ws.Cells.Item(83, 14).Value = -16059
$ws.Cells.Item(97, 8).Value = 1971.6666
$ws.Cells.Item(97, 10).Value = 2166
$ws.Cells.Item(97, 12).Value = 6498
$ws.Cells.Item(97, 14).Value = -7490
$ws.Cells.Item(137, 8).Value = 9039137
$ws.Cells.Item(137, 9).Value = 528468.9399999999
$ws.Cells.Item(137, 11).Value = 1585406.82
$ws.Cells.Item(137, 13).Value = -1582856.82
$ws.Cells.Item(141, 8).Value = 3868
$ws.Cells.Item(141, 9).Value = 3868
$ws.Cells.Item(141, 11).Value = 11604
$ws.Cells.Item(141, 13).Value = -6424
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16401.55
$ws.Cells.Item(32, 9).Value = 17698.773
$ws.Cells.Item(32, 10).Value = 12324.571
$ws.Cells.Item(32, 11).Value = 17698.773
$ws.Cells.Item(32, 12).Value = 12324.571
$ws.Cells.Item(32, 13).Value = -17411.773
$ws.Cells.Item(32, 14).Value = -12898.571
$ws.Cells.Item(61, 8).Value = 6470.674
$ws.Cells.Item(61, 10).Value = 3720.25
$ws.Cells.Item(61, 12).Value = 3720.25
$ws.Cells.Item(61, 14).Value = -4144.25
$ws.Cells.Item(136, 8).Value = 6470.674
$ws.Cells.Item(136, 10).Value = 3720.25
$ws.Cells.Item(136, 12).Value = 11160.75
$ws.Cells.Item(136, 14).Value = -16260.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 442.1
$ws.Cells.Item(80, 9).Value = 744.5
$ws.Cells.Item(80, 10).Value = 366.5
$ws.Cells.Item(80, 11).Value = 744.5
$ws.Cells.Item(80, 12).Value = 366.5
$ws.Cells.Item(80, 13).Value = 253.5
$ws.Cells.Item(80, 14).Value = -2362.5
$ws.Cells.Item(83, 8).Value = 442.1
$ws.Cells.Item(83, 9).Value = 744.5
$ws.Cells.Item(83, 10).Value = 366.5
$ws.Cells.Item(83, 11).Value = 3722.5
$ws.Cells.Item(83, 12).Value = 1832.5
$ws.Cells.Item(83, 13).Value = 1269.5
$ws.Cells.Item(83, 14).Value = -11816.5
$ws.Cells.Item(99, 8).Value = 1737636
$ws.Cells.Item(99, 9).Value = 2977433.2
$ws.Cells.Item(99, 11).Value = 2977433.2
$ws.Cells.Item(99, 13).Value = -2975935.2
$ws.Cells.Item(134, 8).Value = 1629.5883
$ws.Cells.Item(134, 9).Value = 1168.9375
$ws.Cells.Item(134, 11).Value = 3506.8125
$ws.Cells.Item(134, 13).Value = -971.8125
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5713.174
$ws.Cells.Item(31, 10).Value = 7098.4
$ws.Cells.Item(31, 12).Value = 7098.4
$ws.Cells.Item(31, 14).Value = -7688.4
$ws.Cells.Item(34, 8).Value = 5713.174
$ws.Cells.Item(34, 10).Value = 7098.4
$ws.Cells.Item(34, 12).Value = 7098.4
$ws.Cells.Item(34, 14).Value = -7502.4
$ws.Cells.Item(132, 8).Value = 13347395
$ws.Cells.Item(132, 10).Value = 25599.75
$ws.Cells.Item(132, 12).Value = 76799.25
$ws.Cells.Item(132, 14).Value = -81859.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58, 8).Value = 10248.75
$ws.Cells.Item(58, 9).Value = 1999
$ws.Cells.Item(58, 10).Value = 12998.667
$ws.Cells.Item(58, 11).Value = 5997
$ws.Cells.Item(58, 12).Value = 38996.001
$ws.Cells.Item(58, 13).Value = -5869
$ws.Cells.Item(58, 14).Value = -39252.001
$ws.Cells.Item(113, 8).Value = 710.3333
$ws.Cells.Item(113, 9).Value = 682
$ws.Cells.Item(113, 10).Value = 735.125
$ws.Cells.Item(113, 11).Value = 2046
$ws.Cells.Item(113, 12).Value = 2205.375
$ws.Cells.Item(113, 13).Value = 124
$ws.Cells.Item(113, 14).Value = -6545.375
$ws.Cells.Item(139, 8).Value = 1878.75
$ws.Cells.Item(139, 9).Value = 1504.4286
$ws.Cells.Item(139, 11).Value = 4513.2858
$ws.Cells.Item(139, 13).Value = 626.7142000000003
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 3698
$ws.Cells.Item(33, 9).Value = 2000
$ws.Cells.Item(33, 10).Value = 4830
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 12).Value = 4830
$ws.Cells.Item(33, 13).Value = -1748
$ws.Cells.Item(33, 14).Value = -5334
$ws.Cells.Item(70, 8).Value = 5686963.5
$ws.Cells.Item(70, 9).Value = 11367386
$ws.Cells.Item(70, 11).Value = 11367386
$ws.Cells.Item(70, 13).Value = -11367116
$ws.Cells.Item(73, 8).Value = 5686963.5
$ws.Cells.Item(73, 9).Value = 11367386
$ws.Cells.Item(73, 11).Value = 11367386
$ws.Cells.Item(73, 13).Value = -11366450
$ws.Cells.Item(97, 8).Value = 734
$ws.Cells.Item(97, 9).Value = 581.9167
$ws.Cells.Item(97, 10).Value = 1038.1666
$ws.Cells.Item(97, 11).Value = 581.9167
$ws.Cells.Item(97, 12).Value = 1038.1666
$ws.Cells.Item(97, 13).Value = -85.91669999999999
$ws.Cells.Item(97, 14).Value = -2030.1666
$ws.Cells.Item(102, 8).Value = 19238740
$ws.Cells.Item(102, 9).Value = 23817106
$ws.Cells.Item(102, 11).Value = 23817106
$ws.Cells.Item(102, 13).Value = -23815484
$ws.Cells.Item(132, 8).Value = 62766.44
$ws.Cells.Item(132, 9).Value = 86061.625
$ws.Cells.Item(132, 11).Value = 258184.875
$ws.Cells.Item(132, 13).Value = -255654.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1082.6154
$ws.Cells.Item(22, 10).Value = 1479.8
$ws.Cells.Item(22, 12).Value = 1479.8
$ws.Cells.Item(22, 14).Value = -2069.8
$ws.Cells.Item(27, 8).Value = 1082.6154
$ws.Cells.Item(27, 10).Value = 1479.8
$ws.Cells.Item(27, 12).Value = 1479.8
$ws.Cells.Item(27, 14).Value = -1693.8
$ws.Cells.Item(40, 8).Value = 47626476
$ws.Cells.Item(40, 10).Value = 47626476
$ws.Cells.Item(40, 12).Value = 47626476
$ws.Cells.Item(40, 14).Value = -47626748
$ws.Cells.Item(82, 8).Value = 3907242.8
$ws.Cells.Item(82, 9).Value = 6250719
$ws.Cells.Item(82, 11).Value = 6250719
$ws.Cells.Item(82, 13).Value = -6250358
$ws.Cells.Item(85, 8).Value = 3907242.8
$ws.Cells.Item(85, 9).Value = 6250719
$ws.Cells.Item(85, 11).Value = 6250719
$ws.Cells.Item(85, 13).Value = -6249471
$ws.Cells.Item(132, 8).Value = 5045.3716
$ws.Cells.Item(132, 9).Value = 4408.1333
$ws.Cells.Item(132, 10).Value = 5523.3
$ws.Cells.Item(132, 11).Value = 13224.3999
$ws.Cells.Item(132, 12).Value = 16569.9
$ws.Cells.Item(132, 13).Value = -10694.3999
$ws.Cells.Item(132, 14).Value = -21629.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 20005750
$ws.Cells.Item(132, 9).Value = 1691.2941
$ws.Cells.Item(132, 10).Value = 62514376
$ws.Cells.Item(132, 11).Value = 5073.8823
$ws.Cells.Item(132, 12).Value = 187543128
$ws.Cells.Item(132, 13).Value = -2543.8823
$ws.Cells.Item(132, 14).Value = -187548188
$ws.Cells.Item(136, 8).Value = 9061.4375
$ws.Cells.Item(136, 9).Value = 3342.1904
$ws.Cells.Item(136, 11).Value = 10026.5712
$ws.Cells.Item(136, 13).Value = -7476.5712
